$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1. Title / meta-description headline text.
#    Appears twice in the doc: the top Heading1 and a bold run repeated near the
#    bottom of the page (just above the meta description) - wdReplaceAll handles both.
Replace-Text "Play Mad Cars for Free - Review & Game Breakdown" `
             "Play Mad Cars Slot Free - Unique Post-Apocalyptic Game"

# 2. "What we like" bullet list
Replace-Text "Affordable minimum bet at € 0.10" `
             "Unique post-apocalyptic theme"

Replace-Text "Includes special symbols for high-value winnings" `
             "Elongated game grid with 50 fixed pay lines"

Replace-Text "Offers a Mad Bonus Feature as a high-paying car race" `
             "Special symbols for high-value winnings"

# 3. "What we don't like" bullet list
Replace-Text "No Jokers or Scatters for free spins" `
             "No free spin feature"

Replace-Text "Medium volatility may not appeal to those seeking high-risk games" `
             "Limited number of bonus features"

# 4. Meta description (italic run at the very end of the document)
Replace-Text "Read our expert review of Mad Cars, a post-apocalyptic online slot game. Play for free and discover the game's winning potential, bonus features, and more." `
             "Read our review of Mad Cars, a unique post-apocalyptic themed slot game. Play for free on all devices."

Write-Host "Replacements complete"
